# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-09 16:15:32
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Recorded By" email list re-ordered (same people)
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"

# Row 3: a second recorder was added and the attendance count was updated accordingly
$ws.Range("G3").Value = "eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H3").Value = "38/251"

# Average Attendance % figures recalculated after the row-3 update.
# These cells store a literal percentage STRING (not a numeric percent),
# so route the write through a text-formatted helper cell + PasteSpecial
# (values only) to avoid Excel's automatic "24.6%" -> 0.246 numeric coercion
# while still keeping the destination cell's existing style/number format.
$helper = $ws.Range("U1")
$helper.NumberFormat = "@"
$helper.Value = "24.6%"
$helper.Copy()

$ws.Range("L10").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("S15").PasteSpecial(-4163)  # xlPasteValues

$helper.Clear()
$excel.CutCopyMode = $false

# Row 28: "Recorded By" email list re-ordered (same people)
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
